# Update row 8 metrics in metricas_recorrencia_anual sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 1318
$ws.Range("D8").Value = 209
$ws.Range("E8").Value = 1109
$ws.Range("F8").Value = 8.572600492206726
$ws.Range("G8").Value = 84.14264036418815
$ws.Range("H8").Value = 15.85735963581184
